# Applies the "Doing Updates for Financials" commit: inserts a new fiscal-year
# column at column D of the RRTS sheet (shifting the previous D:K data to
# E:L), copies the number/date formatting from the (new) neighboring column E
# into the freshly inserted column D, and then populates D with the new
# year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D; everything that used to be in
# D:K now lives in E:L.
$ws.Columns("D").Insert()

# The freshly inserted column starts out with the generic default width;
# give it the same width as its neighbor (which used to be the first data
# column) so it keeps lining up with the rest of the data block.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Carry over the row-specific number formatting (date format for the header
# rows, #,##0 format for the data rows) from column E into the newly
# inserted column D, restricted to the three data blocks that actually hold
# content so we don't touch untouched rows (5, 6, 36, 37, 78, 79, ...).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the new column D with the newest fiscal year's data (period
# ending 2018-12-31, serial 43465) for the Income Statement, Balance Sheet
# and Cash Flow Statement sections.

# Income Statement (rows 7-35)
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2216100
$ws.Range("D9").Value = 1518400
$ws.Range("D10").Value = 697700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 6200
$ws.Range("D15").Value = 42800
$ws.Range("D17").Value = 2274600
$ws.Range("D18").Value = -58500
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = -15000
$ws.Range("D22").Value = 116900
$ws.Range("D23").Value = -175400
$ws.Range("D24").Value = -9800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -165600
$ws.Range("D27").Value = -165600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = -165600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -165600

# Balance Sheet (rows 38-77)
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 11200
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 278800
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 61100
$ws.Range("D46").Value = 351000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 188700
$ws.Range("D49").Value = 307400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 6400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 853500
$ws.Range("D57").Value = 160200
$ws.Range("D58").Value = 26400
$ws.Range("D59").Value = 110900
$ws.Range("D60").Value = 297600
$ws.Range("D61").Value = 596200
$ws.Range("D62").Value = 11800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 905600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -457400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -52200
$ws.Range("D77").Value = 0

# Cash Flow Statement (rows 80-102)
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -165600
$ws.Range("D83").Value = 43500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 5600
$ws.Range("D91").Value = -25500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -22700
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 2600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -14500
